$d = $word.ActiveDocument

# Texts of the paragraphs that must be removed entirely (including the
# blank paragraph that follows the "LOB1037" requirement line and the
# trailing "Ver no Jupiter ..." / "(c) 2020 ..." footer lines).
$targets = @(
    "Ver no Jupiter Salvar em pdf Salvar em docx",
    [char]0x00A9 + " 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. Original theme under Creative Commons Attribution"
)

foreach ($target in $targets) {
    $paras = $d.Paragraphs
    for ($i = $paras.Count; $i -ge 1; $i--) {
        $para = $paras.Item($i)
        $text = $para.Range.Text
        $trimmed = $text.TrimEnd([char]13, [char]7)
        if ($trimmed -eq $target) {
            $para.Range.Delete()
            break
        }
    }
}

# Remove the now-orphaned blank paragraph that used to separate the
# "LOB1037" requirement line from the "Ver no Jupiter ..." footer block.
$marker = "LOB1037: " + [char]0x00C0 + "lgebra Linear (Requisito fraco)"
$paras = $d.Paragraphs
for ($i = 1; $i -le $paras.Count; $i++) {
    $para = $paras.Item($i)
    $text = $para.Range.Text
    $trimmed = $text.TrimEnd([char]13, [char]7)
    if ($trimmed -eq $marker) {
        $nextPara = $paras.Item($i + 1)
        $nextTrimmed = $nextPara.Range.Text.TrimEnd([char]13, [char]7)
        if ($nextTrimmed -eq "") {
            $nextPara.Range.Delete()
        }
        break
    }
}
